# 10Th - MB for single stock and added new group
# Insert 3 new columns before the old column E (dates), shifting the old
# "Jun_10" column from E to H, then populate the 3 newly-inserted date
# columns (Jun_26, Jun_26, Jun_17 -> actually Jun_17/Jun_15/Jun_13 stay,
# new ones are Jun_27/Jun_26) and fill the new cells with "UN" like the
# rest of the table. Finally append two new rating rows (Benchmark,
# Evercore ISI) as a new broker group.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank columns at E:G - this shifts the existing column E
# (and everything to its right) three columns over, to H.
$ws.Range("E1:G1").EntireColumn.Insert()

# Match the width of the other data columns (C, D, H all render at 8.0).
$ws.Range("E1:G1").EntireColumn.ColumnWidth = 7.14

# --- Header row -------------------------------------------------------
# Old header (before insert) was: B1=Jun_17 C1=Jun_15 D1=Jun_13 E1=Jun_10
# After the column insert, the old E1 (Jun_10) now lives in H1, and the
# rest slid right: B1=Jun_17 C1=Jun_15 D1=Jun_13, new E1:G1 are blank.
# Target final header: B1=Jun_27 C1=Jun_26 D1=Jun_26 E1=Jun_17 F1=Jun_15
# G1=Jun_13 H1=Jun_10 (H1 already correct from the shift).
$ws.Range("E1").Value = $ws.Range("B1").Value()
$ws.Range("F1").Value = $ws.Range("C1").Value()
$ws.Range("G1").Value = $ws.Range("D1").Value()
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# --- Data rows ----------------------------------------------------------
# Fill the three newly inserted columns (E:G) with "UN" for every data
# row, matching the existing "not updated" placeholder used elsewhere in
# the table.
$lastRow = $ws.UsedRange.Rows.Count()
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 5).Value = "UN"
    $ws.Cells.Item($r, 6).Value = "UN"
    $ws.Cells.Item($r, 7).Value = "UN"
}

# --- New broker group rows ----------------------------------------------
$newRow1 = $lastRow + 1
$newRow2 = $lastRow + 2
$ws.Cells.Item($newRow1, 1).Value = "Benchmark"
$ws.Cells.Item($newRow1, 2).Value = "UN"
$ws.Cells.Item($newRow1, 3).Value = "UN"
$ws.Cells.Item($newRow1, 4).Value = "UN"

$ws.Cells.Item($newRow2, 1).Value = "Evercore ISI"
$ws.Cells.Item($newRow2, 2).Value = "UN"
$ws.Cells.Item($newRow2, 3).Value = "UN"
$ws.Cells.Item($newRow2, 4).Value = "UN"
